$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 3.949199438095093
$ws.Range("B1").Value = 2.976226091384888
$ws.Range("C1").Value = 2.727855205535889
$ws.Range("D1").Value = 3.474883079528809
$ws.Range("E1").Value = 4.893055438995361
